$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column header labels to match Canvas naming
$ws.Range("B2:B5").Value = "Exams Final Score"
$ws.Range("C2:C5").Value = "Projects Final Score"

# Reflect the final selection state (entire row 1 selected)
$ws.Rows("1:1").Select()
